$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 2).Value = 0.08254928555990659
$ws.Cells.Item(2, 4).Value = 0.0199032646041033
$ws.Cells.Item(2, 5).Value = 0.1454144946159843
$ws.Cells.Item(2, 6).Value = 0.6222316907933063
$ws.Cells.Item(2, 7).Value = 0.002406805526646769
$ws.Cells.Item(2, 9).Value = 0.4391135274961098
$ws.Cells.Item(2, 11).Value = 0.5193990263817057
$ws.Cells.Item(2, 13).Value = 0.2658377735403263
$ws.Cells.Item(2, 14).Value = 1.294601807314994
$ws.Cells.Item(2, 15).Value = 2.076451734950098

$ws.Cells.Item(3, 2).Value = 0.07309802175362279
$ws.Cells.Item(3, 4).Value = 0.01835804970919241
$ws.Cells.Item(3, 5).Value = 0.1372387361099925
$ws.Cells.Item(3, 6).Value = 0.6158934162451928
$ws.Cells.Item(3, 7).Value = 0.002409375702234238
$ws.Cells.Item(3, 9).Value = 0.4449566295214264
$ws.Cells.Item(3, 11).Value = 0.4529371040189858
$ws.Cells.Item(3, 13).Value = 0.2371435048373272
$ws.Cells.Item(3, 14).Value = 1.308859315523067
$ws.Cells.Item(3, 15).Value = 2.068815343735537

$ws.Cells.Item(4, 2).Value = 0.06729474701933214
$ws.Cells.Item(4, 4).Value = 0.01740150040771482
$ws.Cells.Item(4, 5).Value = 0.1323325121669185
$ws.Cells.Item(4, 6).Value = 0.6124254215324356
$ws.Cells.Item(4, 7).Value = 0.002411038653411015
$ws.Cells.Item(4, 9).Value = 0.4487616592074026
$ws.Cells.Item(4, 11).Value = 0.4120083196593782
$ws.Cells.Item(4, 13).Value = 0.2195887546826327
$ws.Cells.Item(4, 14).Value = 1.318087914377671
$ws.Cells.Item(4, 15).Value = 2.06550874901157

$ws.Cells.Item(5, 2).Value = 0.06493002385433044
$ws.Cells.Item(5, 4).Value = 0.01700976431160228
$ws.Cells.Item(5, 5).Value = 0.130361566957184
$ws.Cells.Item(5, 6).Value = 0.6111186618267581
$ws.Cells.Item(5, 7).Value = 0.002411737720820428
$ws.Cells.Item(5, 9).Value = 0.4503668850378766
$ws.Cells.Item(5, 11).Value = 0.3952997952821988
$ws.Cells.Item(5, 13).Value = 0.2124510300203468
$ws.Cells.Item(5, 14).Value = 1.3219679468418
$ws.Cells.Item(5, 15).Value = 2.064508502053911

$ws.Cells.Item(6, 2).Value = 0.06453737879508026
$ws.Cells.Item(6, 4).Value = 0.01694460064572922
$ws.Cells.Item(6, 5).Value = 0.1300360008288095
$ws.Cells.Item(6, 6).Value = 0.6109081041180744
$ws.Cells.Item(6, 7).Value = 0.002411855094963509
$ws.Cells.Item(6, 9).Value = 0.4506367318068865
$ws.Cells.Item(6, 11).Value = 0.392523583271327
$ws.Cells.Item(6, 13).Value = 0.2112667806331956
$ws.Cells.Item(6, 14).Value = 1.322619429078109
$ws.Cells.Item(6, 15).Value = 2.064363372883406

$ws.Cells.Item(7, 2).Value = 0.06726285463091131
$ws.Cells.Item(7, 4).Value = 0.01739622512131689
$ws.Cells.Item(7, 5).Value = 0.1323058166827664
$ws.Cells.Item(7, 6).Value = 0.6124073670954644
$ws.Cells.Item(7, 7).Value = 0.002411047994682105
$ws.Cells.Item(7, 9).Value = 0.4487830865965181
$ws.Cells.Item(7, 11).Value = 0.4117831020483891
$ws.Cells.Item(7, 13).Value = 0.2194924281680457
$ws.Cells.Item(7, 14).Value = 1.318139758791704
$ws.Cells.Item(7, 15).Value = 2.065493853921595

$ws.Cells.Item(8, 2).Value = 0.07929063935668523
$ws.Cells.Item(8, 4).Value = 0.0193721023189326
$ws.Cells.Item(8, 5).Value = 0.1425717566875662
$ws.Cells.Item(8, 6).Value = 0.6199582595897013
$ws.Cells.Item(8, 7).Value = 0.002407674150678762
$ws.Cells.Item(8, 9).Value = 0.4410831162888957
$ws.Cells.Item(8, 11).Value = 0.496508440066691
$ws.Cells.Item(8, 13).Value = 0.2559307552359655
$ws.Cells.Item(8, 14).Value = 1.299419331246153
$ws.Cells.Item(8, 15).Value = 2.073531620193478

$ws.Cells.Item(9, 2).Value = 0.1028682196504462
$ws.Cells.Item(9, 4).Value = 0.02318424736155578
$ws.Cells.Item(9, 5).Value = 0.1636164817211991
$ws.Cells.Item(9, 6).Value = 0.6381327759145634
$ws.Cells.Item(9, 7).Value = 0.002401728387462603
$ws.Cells.Item(9, 9).Value = 0.4277081263449283
$ws.Cells.Item(9, 11).Value = 0.661672947783245
$ws.Cells.Item(9, 13).Value = 0.3278961086860051
$ws.Cells.Item(9, 14).Value = 1.266472399051413
$ws.Cells.Item(9, 15).Value = 2.100280233298975

$ws.Cells.Item(10, 2).Value = 0.1201767269206471
$ws.Cells.Item(10, 4).Value = 0.02594606777574171
$ws.Cells.Item(10, 5).Value = 0.1796515033590111
$ws.Cells.Item(10, 6).Value = 0.6535479805105524
$ws.Cells.Item(10, 7).Value = 0.002397764575459562
$ws.Cells.Item(10, 9).Value = 0.4189327964422596
$ws.Cells.Item(10, 11).Value = 0.7824023906995308
$ws.Cells.Item(10, 13).Value = 0.3810920505553526
$ws.Cells.Item(10, 14).Value = 1.244559262900559
$ws.Cells.Item(10, 15).Value = 2.126663210091237

$ws.Cells.Item(11, 2).Value = 0.1280460803206722
$ws.Cells.Item(11, 4).Value = 0.02719386865184958
$ws.Cells.Item(11, 5).Value = 0.1870744668818318
$ws.Cells.Item(11, 6).Value = 0.6610108080093937
$ws.Cells.Item(11, 7).Value = 0.002396048301764848
$ws.Cells.Item(11, 9).Value = 0.4151689179646247
$ws.Cells.Item(11, 11).Value = 0.8371881589744419
$ws.Cells.Item(11, 13).Value = 0.4053653208289916
$ws.Cells.Item(11, 14).Value = 1.235088010160357
$ws.Cells.Item(11, 15).Value = 2.140134529557571

$ws.Cells.Item(12, 2).Value = 0.1310251900639656
$ws.Cells.Item(12, 4).Value = 0.02766512725912662
$ws.Cells.Item(12, 5).Value = 0.1899040978761732
$ws.Cells.Item(12, 6).Value = 0.6639016766469581
$ws.Cells.Item(12, 7).Value = 0.002395410821181873
$ws.Cells.Item(12, 9).Value = 0.4137764386456642
$ws.Cells.Item(12, 11).Value = 0.8579141981907128
$ws.Cells.Item(12, 13).Value = 0.4145677554467539
$ws.Cells.Item(12, 14).Value = 1.231572990528459
$ws.Cells.Item(12, 15).Value = 2.145447582395548

$ws.Cells.Item(13, 2).Value = 0.1303836264120406
$ws.Cells.Item(13, 4).Value = 0.02756368956521271
$ws.Cells.Item(13, 5).Value = 0.1892938501540442
$ws.Cells.Item(13, 6).Value = 0.6632761902924216
$ws.Cells.Item(13, 7).Value = 0.002395547561793615
$ws.Cells.Item(13, 9).Value = 0.4140748739211837
$ws.Cells.Item(13, 11).Value = 0.8534513830235539
$ws.Cells.Item(13, 13).Value = 0.4125853716402048
$ws.Cells.Item(13, 14).Value = 1.232326830394072
$ws.Cells.Item(13, 15).Value = 2.144293896569877

$ws.Cells.Item(14, 2).Value = 0.1282911916077438
$ws.Cells.Item(14, 4).Value = 0.02723266470553654
$ws.Cells.Item(14, 5).Value = 0.1873068859144809
$ws.Cells.Item(14, 6).Value = 0.6612473411250477
$ws.Cells.Item(14, 7).Value = 0.002395995606891889
$ws.Cells.Item(14, 9).Value = 0.415053699954667
$ws.Cells.Item(14, 11).Value = 0.8388937108301207
$ws.Cells.Item(14, 13).Value = 0.4061221966938149
$ws.Cells.Item(14, 14).Value = 1.234797393352228
$ws.Cells.Item(14, 15).Value = 2.140567391212471

$ws.Cells.Item(15, 2).Value = 0.1270093986441339
$ws.Cells.Item(15, 4).Value = 0.02702973804932185
$ws.Cells.Item(15, 5).Value = 0.186092257132934
$ws.Cells.Item(15, 6).Value = 0.6600130616285753
$ws.Cells.Item(15, 7).Value = 0.00239627166503006
$ws.Cells.Item(15, 9).Value = 0.4156575339976172
$ws.Cells.Item(15, 11).Value = 0.8299740700694258
$ws.Cells.Item(15, 13).Value = 0.4021647057300157
$ws.Cells.Item(15, 14).Value = 1.236320002412537
$ws.Cells.Item(15, 15).Value = 2.138312388900374

$ws.Cells.Item(16, 2).Value = 0.119662339802332
$ws.Cells.Item(16, 4).Value = 0.02586434657698078
$ws.Cells.Item(16, 5).Value = 0.1791690017002665
$ws.Cells.Item(16, 6).Value = 0.6530693379585983
$ws.Cells.Item(16, 7).Value = 0.002397878481151876
$ws.Cells.Item(16, 9).Value = 0.4191833660640221
$ws.Cells.Item(16, 11).Value = 0.7788192356463526
$ws.Cells.Item(16, 13).Value = 0.3795072284705157
$ws.Cells.Item(16, 14).Value = 1.245188239508636
$ws.Cells.Item(16, 15).Value = 2.125812436320615

$ws.Cells.Item(17, 2).Value = 0.1151538841580759
$ws.Cells.Item(17, 4).Value = 0.02514720524398939
$ws.Cells.Item(17, 5).Value = 0.1749548976586794
$ws.Cells.Item(17, 6).Value = 0.6489250074590558
$ws.Cells.Item(17, 7).Value = 0.002398886422131567
$ws.Cells.Item(17, 9).Value = 0.4214047797243001
$ws.Cells.Item(17, 11).Value = 0.7474023579158029
$ws.Cells.Item(17, 13).Value = 0.3656266027816457
$ws.Cells.Item(17, 14).Value = 1.250755980377615
$ws.Cells.Item(17, 15).Value = 2.118520821669847

$ws.Cells.Item(18, 2).Value = 0.1125603408842295
$ws.Cells.Item(18, 4).Value = 0.02473391981494188
$ws.Cells.Item(18, 5).Value = 0.1725431425999275
$ws.Cells.Item(18, 6).Value = 0.6465836838482772
$ws.Cells.Item(18, 7).Value = 0.002399474344560895
$ws.Cells.Item(18, 9).Value = 0.4227039430105233
$ws.Cells.Item(18, 11).Value = 0.7293195723568431
$ws.Cells.Item(18, 13).Value = 0.3576498147207658
$ws.Cells.Item(18, 14).Value = 1.254005187936901
$ws.Cells.Item(18, 15).Value = 2.114465174829036

$ws.Cells.Item(19, 2).Value = 0.1116821497444676
$ws.Cells.Item(19, 4).Value = 0.02459385101970213
$ws.Cells.Item(19, 5).Value = 0.1717286325686374
$ws.Cells.Item(19, 6).Value = 0.6457982287904258
$ws.Cells.Item(19, 7).Value = 0.002399674811802071
$ws.Cells.Item(19, 9).Value = 0.4231475032476766
$ws.Cells.Item(19, 11).Value = 0.7231949066449204
$ws.Cells.Item(19, 13).Value = 0.3549502100208031
$ws.Cells.Item(19, 14).Value = 1.255113348597344
$ws.Cells.Item(19, 15).Value = 2.113115741920325

$ws.Cells.Item(20, 2).Value = 0.1156338597922542
$ws.Cells.Item(20, 4).Value = 0.02522362962664459
$ws.Cells.Item(20, 5).Value = 0.1754022436933482
$ws.Cells.Item(20, 6).Value = 0.6493617910679603
$ws.Cells.Item(20, 7).Value = 0.002398778278821301
$ws.Cells.Item(20, 9).Value = 0.4211660847330094
$ws.Cells.Item(20, 11).Value = 0.7507480508177764
$ws.Cells.Item(20, 13).Value = 0.3671034965448214
$ws.Cells.Item(20, 14).Value = 1.25015844188642
$ws.Cells.Item(20, 15).Value = 2.119282710120672

$ws.Cells.Item(21, 2).Value = 0.1289058150360063
$ws.Cells.Item(21, 4).Value = 0.02732992903037257
$ws.Cells.Item(21, 5).Value = 0.1878899956988889
$ws.Cells.Item(21, 6).Value = 0.661841502331896
$ws.Cells.Item(21, 7).Value = 0.002395863668138682
$ws.Cells.Item(21, 9).Value = 0.4147653041059947
$ws.Cells.Item(21, 11).Value = 0.843170204889617
$ws.Cells.Item(21, 13).Value = 0.4080202980911309
$ws.Cells.Item(21, 14).Value = 1.234069787186971
$ws.Cells.Item(21, 15).Value = 2.141656205700485

$ws.Cells.Item(22, 2).Value = 0.137574788468811
$ws.Cells.Item(22, 4).Value = 0.02869918205917799
$ws.Cells.Item(22, 5).Value = 0.1961606605966892
$ws.Cells.Item(22, 6).Value = 0.6703757995411621
$ws.Cells.Item(22, 7).Value = 0.002394031255965409
$ws.Cells.Item(22, 9).Value = 0.4107733341792876
$ws.Cells.Item(22, 11).Value = 0.9034556602942985
$ws.Cells.Item(22, 13).Value = 0.4348240129544507
$ws.Cells.Item(22, 14).Value = 1.223971925905161
$ws.Cells.Item(22, 15).Value = 2.157513058494885

$ws.Cells.Item(23, 2).Value = 0.1329485240627122
$ws.Cells.Item(23, 4).Value = 0.02796906534646126
$ws.Cells.Item(23, 5).Value = 0.1917363834128025
$ws.Cells.Item(23, 6).Value = 0.6657862586661736
$ws.Cells.Item(23, 7).Value = 0.002395002638595927
$ws.Cells.Item(23, 9).Value = 0.4128864106147141
$ws.Cells.Item(23, 11).Value = 0.8712912106515773
$ws.Cells.Item(23, 13).Value = 0.4205126700581019
$ws.Cells.Item(23, 14).Value = 1.229323171407064
$ws.Cells.Item(23, 15).Value = 2.148936859483172

$ws.Cells.Item(24, 2).Value = 0.1154168676881682
$ws.Cells.Item(24, 4).Value = 0.02518908124916663
$ws.Cells.Item(24, 5).Value = 0.1751999643396118
$ws.Cells.Item(24, 6).Value = 0.6491641925353093
$ws.Cells.Item(24, 7).Value = 0.002398827143964706
$ws.Cells.Item(24, 9).Value = 0.4212739300779216
$ws.Cells.Item(24, 11).Value = 0.7492355278031653
$ws.Cells.Item(24, 13).Value = 0.3664357823864179
$ws.Cells.Item(24, 14).Value = 1.250428438776904
$ws.Cells.Item(24, 15).Value = 2.118937835502436

$ws.Cells.Item(25, 2).Value = 0.09649169269610525
$ws.Cells.Item(25, 4).Value = 0.02215974317834934
$ws.Cells.Item(25, 5).Value = 0.1578238289863449
$ws.Cells.Item(25, 6).Value = 0.6328546044304773
$ws.Cells.Item(25, 7).Value = 0.002403265536989977
$ws.Cells.Item(25, 9).Value = 0.4311417936037838
$ws.Cells.Item(25, 11).Value = 0.6170985765418777
$ws.Cells.Item(25, 13).Value = 0.3083717285726237
$ws.Cells.Item(25, 14).Value = 1.274982648397078
$ws.Cells.Item(25, 15).Value = 2.091864418967504
